$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 21.137664
$ws.Range("H2").Value = 63.412992
$ws.Range("I2").Value = 0.01636081711441431
$ws.Range("J2").Value = 0.01636081711441431
$ws.Range("M2").Value = 12.67919733333333
$ws.Range("N2").Value = 38.037592
$ws.Range("O2").Value = 0.9871416146107245
$ws.Range("P2").Value = 0.9871416146107247
$ws.Range("Q2").Value = 268.008613021696
$ws.Range("R2").Value = 2412.077517195264
$ws.Range("S2").Value = 0.01615044342267371
$ws.Range("T2").Value = 0.01615044342267372

# Row 3
$ws.Range("G3").Value = 21.137664
$ws.Range("H3").Value = 63.412992
$ws.Range("I3").Value = 0.01636081711441431
$ws.Range("J3").Value = 0.01636081711441431
$ws.Range("M3").Value = 0.1651576666666667
$ws.Range("N3").Value = 0.495473
$ws.Range("O3").Value = 0.01285838538927542
$ws.Range("P3").Value = 0.01285838538927542
$ws.Range("Q3").Value = 3.491047265024001
$ws.Range("R3").Value = 31.419425385216
$ws.Range("S3").Value = 0.0002103736917405921
$ws.Range("T3").Value = 0.0002103736917405921

# Row 4
$ws.Range("I4").Value = 0.8979526429041496
$ws.Range("J4").Value = 0.8979526429041496
$ws.Range("M4").Value = 12.67919733333333
$ws.Range("N4").Value = 38.037592
$ws.Range("O4").Value = 0.9871416146107245
$ws.Range("P4").Value = 0.9871416146107247
$ws.Range("Q4").Value = 14709.47573711832
$ws.Range("R4").Value = 132385.2816340649
$ws.Range("S4").Value = 0.8864064217603695
$ws.Range("T4").Value = 0.8864064217603697

# Row 5
$ws.Range("I5").Value = 0.8979526429041496
$ws.Range("J5").Value = 0.8979526429041496
$ws.Range("M5").Value = 0.1651576666666667
$ws.Range("N5").Value = 0.495473
$ws.Range("O5").Value = 0.01285838538927542
$ws.Range("P5").Value = 0.01285838538927542
$ws.Range("S5").Value = 0.01154622114377996
$ws.Range("T5").Value = 0.01154622114377996

# Row 6
$ws.Range("G6").Value = 109.7535913333333
$ws.Range("H6").Value = 329.260774
$ws.Range("I6").Value = 0.08495065658413503
$ws.Range("J6").Value = 0.08495065658413503
$ws.Range("M6").Value = 12.67919733333333
$ws.Range("N6").Value = 38.037592
$ws.Range("O6").Value = 0.9871416146107245
$ws.Range("P6").Value = 0.9871416146107247
$ws.Range("Q6").Value = 1391.587442557357
$ws.Range("R6").Value = 12524.28698301621
$ws.Range("S6").Value = 0.08385832830270423
$ws.Range("T6").Value = 0.08385832830270425

# Row 7
$ws.Range("G7").Value = 109.7535913333333
$ws.Range("H7").Value = 329.260774
$ws.Range("I7").Value = 0.08495065658413503
$ws.Range("J7").Value = 0.08495065658413503
$ws.Range("M7").Value = 0.1651576666666667
$ws.Range("N7").Value = 0.495473
$ws.Range("O7").Value = 0.01285838538927542
$ws.Range("P7").Value = 0.01285838538927542
$ws.Range("Q7").Value = 18.12664705290022
$ws.Range("R7").Value = 163.139823476102
$ws.Range("S7").Value = 0.001092328281430795
$ws.Range("T7").Value = 0.001092328281430795

# Row 8
$ws.Range("G8").Value = 0.9507383333333334
$ws.Range("H8").Value = 2.852215
$ws.Range("I8").Value = 0.0007358833973011272
$ws.Range("J8").Value = 0.0007358833973011273
$ws.Range("M8").Value = 12.67919733333333
$ws.Range("N8").Value = 38.037592
$ws.Range("O8").Value = 0.9871416146107245
$ws.Range("P8").Value = 0.9871416146107247
$ws.Range("Q8").Value = 12.05459894069778
$ws.Range("R8").Value = 108.49139046628
$ws.Range("S8").Value = 0.0007264211249770599
$ws.Range("T8").Value = 0.0007264211249770603

# Row 9
$ws.Range("G9").Value = 0.9507383333333334
$ws.Range("H9").Value = 2.852215
$ws.Range("I9").Value = 0.0007358833973011272
$ws.Range("J9").Value = 0.0007358833973011273
$ws.Range("M9").Value = 0.1651576666666667
$ws.Range("N9").Value = 0.495473
$ws.Range("O9").Value = 0.01285838538927542
$ws.Range("P9").Value = 0.01285838538927542
$ws.Range("Q9").Value = 0.1570217247438889
$ws.Range("R9").Value = 1.413195522695
$ws.Range("S9").Value = 0.000009462272324067171
$ws.Range("T9").Value = 0.000009462272324067173
